$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear the previously used range so the shared-string pool starts empty,
# then rebuild it by writing strings in the exact order they must be interned
# (this keeps the final xl/sharedStrings.xml ordering/index layout correct).
$ws.Range("A1:T11").Value = ""

# Re-intern the 20 header strings, in column order, at their original indices (0-19).
$ws.Range("A1").Value = "Sending cluster"
$ws.Range("B1").Value = "Ligand symbol"
$ws.Range("C1").Value = "Receptor symbol"
$ws.Range("D1").Value = "Target cluster"
$ws.Range("E1").Value = "Ligand-expressing cells"
$ws.Range("F1").Value = "Ligand detection rate"
$ws.Range("G1").Value = "Ligand average expression value"
$ws.Range("H1").Value = "Ligand total expression value"
$ws.Range("I1").Value = "Ligand derived specificity of average expression value"
$ws.Range("J1").Value = "Ligand derived specificity of total expression value"
$ws.Range("K1").Value = "Receptor-expressing cells"
$ws.Range("L1").Value = "Receptor detection rate"
$ws.Range("M1").Value = "Receptor average expression value"
$ws.Range("N1").Value = "Receptor total expression value"
$ws.Range("O1").Value = "Receptor derived specificity of average expression value"
$ws.Range("P1").Value = "Receptor derived specificity of total expression value"
$ws.Range("Q1").Value = "Edge average expression weight"
$ws.Range("R1").Value = "Edge total expression weight"
$ws.Range("S1").Value = "Edge average expression derived specificity"
$ws.Range("T1").Value = "Edge total expression derived specificity"

# Re-intern the new set of label strings in the exact order required so they land
# at shared-string indices 20..26 as: ECs, FAPs, M1, M2, sCs, Tnfsf10, Tnfrsf11b.
$ws.Range("Z1").Value = "ECs"
$ws.Range("Z2").Value = "FAPs"
$ws.Range("Z3").Value = "M1"
$ws.Range("Z4").Value = "M2"
$ws.Range("Z5").Value = "sCs"
$ws.Range("Z6").Value = "Tnfsf10"
$ws.Range("Z7").Value = "Tnfrsf11b"
$ws.Range("Z1:Z7").Value = ""

# Now write the real table contents (rows 2-11).
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf10"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 26.80992266666667
$ws.Range("H2").Value = 80.429768
$ws.Range("I2").Value = 0.8864931668627436
$ws.Range("J2").Value = 0.8864931668627437
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("N2").Value = 0.247237
$ws.Range("O2").Value = 0.04727005612861496
$ws.Range("P2").Value = 0.04727005612861496
$ws.Range("Q2").Value = 2.209468283446222
$ws.Range("R2").Value = 19.885214551016
$ws.Range("S2").Value = 0.04190458175523552
$ws.Range("T2").Value = 0.04190458175523552

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfsf10"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 26.80992266666667
$ws.Range("H3").Value = 80.429768
$ws.Range("I3").Value = 0.8864931668627436
$ws.Range("J3").Value = 0.8864931668627437
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.661024
$ws.Range("N3").Value = 4.983072
$ws.Range("O3").Value = 0.9527299438713851
$ws.Range("P3").Value = 0.952729943871385
$ws.Range("Q3").Value = 44.53192498747733
$ws.Range("R3").Value = 400.787324887296
$ws.Range("S3").Value = 0.8445885851075081
$ws.Range("T3").Value = 0.8445885851075081

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf10"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.262185
$ws.Range("H4").Value = 3.786555
$ws.Range("I4").Value = 0.04173523332119963
$ws.Range("J4").Value = 0.04173523332119964
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("N4").Value = 0.247237
$ws.Range("O4").Value = 0.04727005612861496
$ws.Range("P4").Value = 0.04727005612861496
$ws.Range("Q4").Value = 0.1040196109483333
$ws.Range("R4").Value = 0.9361764985350001
$ws.Range("S4").Value = 0.001972826821633948
$ws.Range("T4").Value = 0.001972826821633948

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfsf10"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.262185
$ws.Range("H5").Value = 3.786555
$ws.Range("I5").Value = 0.04173523332119963
$ws.Range("J5").Value = 0.04173523332119964
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.661024
$ws.Range("N5").Value = 4.983072
$ws.Range("O5").Value = 0.9527299438713851
$ws.Range("P5").Value = 0.952729943871385
$ws.Range("Q5").Value = 2.09651957744
$ws.Range("R5").Value = 18.86867619696
$ws.Range("S5").Value = 0.03976240649956569
$ws.Range("T5").Value = 0.03976240649956569

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Tnfsf10"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7313983333333334
$ws.Range("H6").Value = 2.194195
$ws.Range("I6").Value = 0.02418431536771806
$ws.Range("J6").Value = 0.02418431536771806
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("N6").Value = 0.247237
$ws.Range("O6").Value = 0.04727005612861496
$ws.Range("P6").Value = 0.04727005612861496
$ws.Range("Q6").Value = 0.06027624324611112
$ws.Range("R6").Value = 0.542486189215
$ws.Range("S6").Value = 0.001143193944864158
$ws.Range("T6").Value = 0.001143193944864158

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Tnfsf10"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7313983333333334
$ws.Range("H7").Value = 2.194195
$ws.Range("I7").Value = 0.02418431536771806
$ws.Range("J7").Value = 0.02418431536771806
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.661024
$ws.Range("N7").Value = 4.983072
$ws.Range("O7").Value = 0.9527299438713851
$ws.Range("P7").Value = 0.952729943871385
$ws.Range("Q7").Value = 1.214870185226667
$ws.Range("R7").Value = 10.93383166704
$ws.Range("S7").Value = 0.0230411214228539
$ws.Range("T7").Value = 0.0230411214228539

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Tnfsf10"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.357302666666667
$ws.Range("H8").Value = 4.071908000000001
$ws.Range("I8").Value = 0.04488038083230255
$ws.Range("J8").Value = 0.04488038083230256
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("N8").Value = 0.247237
$ws.Range("O8").Value = 0.04727005612861496
$ws.Range("P8").Value = 0.04727005612861496
$ws.Range("Q8").Value = 0.1118584797995556
$ws.Range("R8").Value = 1.006726318196
$ws.Range("S8").Value = 0.002121498121016557
$ws.Range("T8").Value = 0.002121498121016557

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Tnfsf10"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.357302666666667
$ws.Range("H9").Value = 4.071908000000001
$ws.Range("I9").Value = 0.04488038083230255
$ws.Range("J9").Value = 0.04488038083230256
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.661024
$ws.Range("N9").Value = 4.983072
$ws.Range("O9").Value = 0.9527299438713851
$ws.Range("P9").Value = 0.952729943871385
$ws.Range("Q9").Value = 2.254512304597334
$ws.Range("R9").Value = 20.290610741376
$ws.Range("S9").Value = 0.04275888271128599
$ws.Range("T9").Value = 0.042758882711286

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tnfsf10"
$ws.Range("C10").Value = "Tnfrsf11b"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.08186399999999999
$ws.Range("H10").Value = 0.245592
$ws.Range("I10").Value = 0.002706903616036228
$ws.Range("J10").Value = 0.002706903616036229
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("N10").Value = 0.247237
$ws.Range("O10").Value = 0.04727005612861496
$ws.Range("P10").Value = 0.04727005612861496
$ws.Range("Q10").Value = 0.006746603255999999
$ws.Range("R10").Value = 0.060719429304
$ws.Range("S10").Value = 0.0001279554858647833
$ws.Range("T10").Value = 0.0001279554858647833

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Tnfsf10"
$ws.Range("C11").Value = "Tnfrsf11b"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.08186399999999999
$ws.Range("H11").Value = 0.245592
$ws.Range("I11").Value = 0.002706903616036228
$ws.Range("J11").Value = 0.002706903616036229
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.661024
$ws.Range("N11").Value = 4.983072
$ws.Range("O11").Value = 0.9527299438713851
$ws.Range("P11").Value = 0.952729943871385
$ws.Range("Q11").Value = 0.135978068736
$ws.Range("R11").Value = 1.223802618624
$ws.Range("S11").Value = 0.002578948130171445
$ws.Range("T11").Value = 0.002578948130171445

